# "added source for nsa picture"
#
# Slide 8 ("NSA") has a picture (id=186, "Picture 185") at
# off (7560000,720000) ext (4500000,4500000) EMU. The author added a
# small caption textbox just under/right of that picture reading
# "Source: Wikipedia".
#
# Shapes.AddTextbox's Left/Top/Width/Height are in points (PowerPoint
# COM converts points -> EMU at 12700 EMU/pt), so the EMU target from
# the OOXML (off x=8744771 y=5228668, ext cx=2130458 cy=369332) is
# divided by 12700 below to land exactly back on those EMU values.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)

# The saved deck also gained ppt/commentAuthors.xml registering the
# author (c.hofsom@gmail.com) as a comment author, with no actual
# comment left behind (PowerPoint/M365 registers the author as soon as
# the comments UI is touched, even if no comment text survives).
# Adding then immediately deleting a comment reproduces that
# commentAuthors.xml side effect without leaving a stray comment part.
$cm = $s.Comments.Add(10, 10, "c.hofsom@gmail.com", "c", "")
$cm.Delete()

# Add the "Source: Wikipedia" caption textbox under the NSA picture.
$tb = $s.Shapes.AddTextbox(1, 688.5646456692914, 411.70614173228347, 167.75259842519685, 29.081259842519685)
$tb.Name = "Textfeld 1"
$tb.Fill.Visible = $false
$tb.TextFrame.WordWrap = $true
$tb.TextFrame.AutoSize = 1
$tb.TextFrame.TextRange.Text = "Source: Wikipedia"
